# TC07_Canine_Filter_SamplePatho-Osteosarcoma.xlsx - "startup" sheet
# Commit: updated keyword for case files tab data writing issue + more curated scripts
#
# - CasesTab (row 2) query: add a `Cohort` column to the RETURN clause.
# - FilesTab (row 4) query: drop the trailing `Study Code` column from the RETURN clause.
# - SamplesTab (row 3) query: left untouched.
# - StatQuery column (C2:C4) for all three tabs: replaced with a new, shorter
#   curated summary script (Programs/Studies/Cases/Samples/Case Files/Study Files)
#   instead of the old verbose filter-driven aliquot/file-count script.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
MATCH (samp:sample)-->(c) 
 WHERE samp.specific_sample_pathology IN ["Osteosarcoma"]  
OPTIONAL MATCH (co:cohort)<-[*]-(c)
  WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`,
coalesce(co.cohort_description, '') AS `Cohort`
'@

$newStatQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
 WHERE samp.specific_sample_pathology IN ["Osteosarcoma"]  
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

$filesQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
 MATCH (samp:sample)-->(c) 
 WHERE samp.specific_sample_pathology IN ["Osteosarcoma"]  
WITH DISTINCT f, parent, c, demo, diag, s
RETURN coalesce(f.file_name, '') AS `File Name`, 
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(diag.disease_term,'') AS Diagnosis 
'@

# Row 2 - CasesTab
$ws.Range("B2").Value = $casesQuery
$ws.Range("C2").Value = $newStatQuery

# Row 3 - SamplesTab (query column B3 is unchanged; only StatQuery is refreshed)
$ws.Range("C3").Value = $newStatQuery

# Row 4 - FilesTab
$ws.Range("B4").Value = $filesQuery
$ws.Range("C4").Value = $newStatQuery

# The new StatQuery text is much shorter, so the wrapped-text rows shrink from
# the capped 409.6pt down to their real auto-fit heights.
$ws.Rows.Item(2).RowHeight = 270
$ws.Rows.Item(3).RowHeight = 225
$ws.Rows.Item(4).RowHeight = 210

# Column widths + zoom level also shifted slightly on the re-save.
$ws.Columns.Item(1).ColumnWidth = 10.85546875
$ws.Columns.Item(2).ColumnWidth = 87.7109375
$ws.Columns.Item(3).ColumnWidth = 75.7109375
$ws.Columns.Item(4).ColumnWidth = 70.28515625
$ws.Columns.Item(5).ColumnWidth = 28.5703125

$ws.Activate()
$excel.ActiveWindow.Zoom = 70
